$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 1 ("Peter Klembczyk"):
#    - Insert a new paragraph before it: "Daniel's Edits in red" in red,
#      with the paragraph mark itself also carrying the red run-properties.
#    - Split "Peter Klembczyk" into "Peter " + "Klembczyk" runs wrapped with
#      proofErr spellStart/spellEnd (as Word's spell checker does for the
#      flagged surname).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:rPr>
<w:color w:val="FF0000"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:color w:val="FF0000"/>
</w:rPr>
<w:t>Daniel&#8217;s Edits in red</w:t>
</w:r>
</w:p>
<w:p>
<w:r>
<w:t xml:space="preserve">Peter </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>Klembczyk</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Mark the inline picture's run as "no proofing" (<w:noProof/>), which is
#    what Word stamps on runs holding drawings/pictures.
# ---------------------------------------------------------------------------
$shp = $d.InlineShapes.Item(1)
$shp.Range.NoProofing = 1

# ---------------------------------------------------------------------------
# 3) Add a <w:lastRenderedPageBreak/> before the "However, for it..." run,
#    and append a new red sub-bullet right after it.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "However, for it to be maintained*") {
        $target = $cand
        break
    }
}

$r2 = $target.Range
$xml2 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="5"/>
</w:numPr>
</w:pPr>
<w:r>
<w:lastRenderedPageBreak/>
<w:t>However, for it to be maintained, tests would have to be designed and the results observed from people using it for a semester</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="1"/>
<w:numId w:val="5"/>
</w:numPr>
</w:pPr>
<w:r>
<w:rPr>
<w:color w:val="FF0000"/>
</w:rPr>
<w:t>This could be anything from testing new users in the system to making sure that users who no longer need access (graduate or leave the school) can&#8217;t still use it</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$r2.InsertXML($xml2)
